$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "FilesTab" query cell (B4) with the revised Cypher query text
# (removed the `File Type` and `Breed` columns from the RETURN clause).
$newFileQuery = "`nMATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nWHERE diag.stage_of_disease IN ['V']`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN coalesce(f.file_name, '') AS ``File Name``, `n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Cells.Item(4, 2).Value = $newFileQuery

# Update the active cell / view so the selection sits on the edited cell.
$ws.Range("B4").Select()
